$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burn Down Chart")

try {
    # New shared string entry appears because "Rest API" is written into D16.
    $ws.Range("D16").Value = "Rest API"
    $ws.Range("E16").Value = "Ali Cooper"

    # Update the daily burn-down effort cells (rows 11-16).
    $ws.Range("G11").Value = 1
    $ws.Range("G12").Value = 0.5
    $ws.Range("H13").Value = 0.5
    $ws.Range("F14").Value = 1
    $ws.Range("H14").Value = 5
    $ws.Range("H15").Value = 1
    $ws.Range("I15").Value = 5
    $ws.Range("J15").Value = 0
    $ws.Range("F16").Value = 2
    $ws.Range("G16").Value = 1
    $ws.Range("H16").Value = 2
    $ws.Range("I16").Value = 1

    # Move the active selection from I12 to G12.
    $ws.Range("G12").Select()
}
catch {
    Write-Host "Error applying Sprint3 burndown updates: $_"
}
